# Add more files from Albany Evening News to the FileData sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("FileData")

# --- Fix up columns M/N on existing rows -------------------------------
# In the original sheet, column M held "tags" and column N held
# "description". Going forward, M holds "description" and N holds "tags",
# so the header row and the one existing data row need their M/N values
# swapped.
$ws.Range("M1").Value = "description"
$ws.Range("N1").Value = "tags"

$ws.Range("M3").Value = "Castleton Scouts To Present Movie"
$ws.Range("N3").Value = "whispering wires"

# --- New rows for the newly-added newspaper clippings -------------------
$ws.Range("A4").Value = "albany_evening_news_1936_jun_25_big_moose.pdf"
$ws.Range("B4").Value = "newspapers"
$ws.Range("C4").Value = 1936
$ws.Range("D4").Value = 6
$ws.Range("E4").Value = 25
$ws.Range("F4").Value = "Albany Evening News"
$ws.Range("H4").Value = "X"
$ws.Range("M4").Value = "Bake Appointed to Big Moose Staff"
$ws.Range("N4").Value = "big moose,thomas latham,donald paul"

$ws.Range("A5").Value = "albany_evening_news_1936_sep_25_promotions.pdf"
$ws.Range("B5").Value = "newspapers"
$ws.Range("C5").Value = 1936
$ws.Range("D5").Value = 9
$ws.Range("E5").Value = 25
$ws.Range("F5").Value = "Albany Evening News"
$ws.Range("H5").Value = "X"
$ws.Range("M5").Value = "Boy Scouts Win New Promotions"
$ws.Range("N5").Value = "fort orange council,court of honor,albert bleadow,peter andrew,thomas latham,christian gersch"

$ws.Range("A6").Value = "albany_evening_news_1937_jun_11_camp_sign_up.pdf"
$ws.Range("B6").Value = "newspapers"
$ws.Range("C6").Value = 1937
$ws.Range("D6").Value = 6
$ws.Range("E6").Value = 11
$ws.Range("F6").Value = "Albany Evening News"
$ws.Range("H6").Value = "X"
$ws.Range("M6").Value = "Boy Scout Camp Enrollment Gains"
$ws.Range("N6").Value = "big moose,fort orange council,william grooten,thomas latham"

$ws.Range("A7").Value = "albany_evening_news_1937_jun_29_camp_sign_up.pdf"
$ws.Range("B7").Value = "newspapers"
$ws.Range("C7").Value = 1937
$ws.Range("D7").Value = 6
$ws.Range("E7").Value = 29
$ws.Range("F7").Value = "Albany Evening News"
$ws.Range("H7").Value = "X"
$ws.Range("M7").Value = "Assistant Director Appointed For Camp"
$ws.Range("N7").Value = "fort orange council,camp hawley,george corson,robert lawrence"

# --- Column widths --------------------------------------------------------
# Widen column A (filenames), and M/N (description/tags) so the new,
# longer values fit.
$ws.Columns.Item(1).ColumnWidth = 47.333333333333336
$ws.Columns.Item(13).ColumnWidth = 36.166666666666664
$ws.Columns.Item(14).ColumnWidth = 87.16666666666667

# --- View / selection ------------------------------------------------------
# Scroll back to the left edge and select A9 (first empty row below the
# new data), matching where the author left off editing.
$ws.Range("A9").Select()
